$d = $word.ActiveDocument

# The four paragraphs that contain stray <w:proofErr/> spell-check markers
# (around "WillieHand", "WillieHandTests", "Exercie" and ".Fail"/"(" ).
# Round-tripping a paragraph's Range through WordOpenXML -> InsertXML drops
# those markers (they're not meaningful outside of full-document context)
# and also coalesces any adjacent runs that share identical formatting,
# which is exactly what the target revision does.
$targets = @(7, 8, 13, 14)
foreach ($i in $targets) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $xml = $rng.WordOpenXML
    $rng.InsertXML($xml)
}

# Append the new bullet paragraph at the end of the document.
$wholeRng = $d.Content
$wholeXml = $wholeRng.WordOpenXML
$wholeXml = $wholeXml.Replace('<w:p w14:paraId="00000001" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"/>', '')
$newParagraph = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add MoreWillieHand and make it partial class WillieHand.</w:t></w:r></w:p>'
$sectPrIndex = $wholeXml.IndexOf('<w:sectPr')
$wholeXml = $wholeXml.Substring(0, $sectPrIndex) + $newParagraph + $wholeXml.Substring($sectPrIndex)
$wholeRng.InsertXML($wholeXml)
